$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new session timestamp.
$ws.Name = "牌局記錄0809_191852"

# New hand-history rows (875-899) appended to the log.
$newRowsData = @"
8|1,7|3|3,J|5|閒|藍|藍|藍|紅|紅|紅
1|6,5|8|10,8|7|莊|紅|紅|紅|藍|藍|藍
3|7,6,J|9|Q,1,8|6|莊|藍|藍|藍|紅|紅|紅
5|K,5,J|0|10,K,K|5|閒|紅|紅|紅|藍|藍|藍
2|2,J,10|6|5,1|4|莊|藍|紅|紅|紅|藍|藍
6|4,2|3|8,2,3|3|閒|藍|紅|藍|紅|藍|紅
9|6,3|9|6,3|0|和|紅|紅|藍|藍|藍|紅
4|3,2,9|8|6,4,8|4|莊|紅|紅|藍|藍|藍|紅
1|K,1|8|9,9|7|莊|藍|藍|藍|紅|紅|紅
8|8,J|8|10,8|0|和|紅|紅|紅|藍|藍|藍
2|9,3,Q|0|10,K,Q|2|閒|紅|紅|紅|藍|藍|藍
4|10,4|8|2,6|4|莊|藍|紅|紅|紅|藍|藍
0|3,7,J|0|7,3,J|0|和|藍|紅|藍|紅|藍|紅
3|4,9,K|1|5,8,8|2|閒|藍|紅|藍|紅|藍|紅
1|1,K,K|0|K,K,K|1|閒|紅|紅|藍|藍|藍|紅
9|1,8|5|K,5|4|閒|藍|藍|紅|紅|紅|藍
2|K,Q,2|5|2,3|3|莊|藍|藍|藍|紅|紅|紅
6|J,6|5|7,7,1|1|閒|紅|藍|藍|藍|紅|紅
2|K,Q,2|5|3,J,2|3|莊|紅|藍|紅|藍|紅|藍
7|10,J,7|2|7,5,J|5|閒|藍|藍|紅|紅|紅|藍
0|1,9,K|9|Q,J,9|9|莊|紅|紅|紅|藍|藍|藍
3|8,7,8|5|8,7|2|莊|藍|藍|藍|藍|藍|藍
1|Q,1,K|4|3,1|3|莊|紅|紅|紅|藍|藍|藍
8|8,K|1|7,4|7|閒|紅|紅|紅|藍|藍|藍
5|5,J|9|9,Q|4|莊|藍|紅|紅|紅|藍|藍
"@

$startRow = 875
$lines = $newRowsData -split "`r?`n" | Where-Object { $_.Length -gt 0 }

$rowIndex = $startRow
foreach ($line in $lines) {
    $fields = $line -split '\|'

    $ws.Cells.Item($rowIndex, 1).Value = [int]$fields[0]
    $ws.Cells.Item($rowIndex, 2).Value = $fields[1]
    $ws.Cells.Item($rowIndex, 3).Value = [int]$fields[2]
    $ws.Cells.Item($rowIndex, 4).Value = $fields[3]
    $ws.Cells.Item($rowIndex, 5).Value = [int]$fields[4]
    $ws.Cells.Item($rowIndex, 6).Value = $fields[5]
    $ws.Cells.Item($rowIndex, 7).Value = $fields[6]
    $ws.Cells.Item($rowIndex, 8).Value = $fields[7]
    $ws.Cells.Item($rowIndex, 9).Value = $fields[8]
    $ws.Cells.Item($rowIndex, 10).Value = $fields[9]
    $ws.Cells.Item($rowIndex, 11).Value = $fields[10]
    $ws.Cells.Item($rowIndex, 12).Value = $fields[11]

    $rowIndex++
}

Write-Output "Wrote rows $startRow to $($rowIndex - 1)"
